# The sheet is protected, so it must be unprotected before the cell
# values can be updated, and re-protected afterwards to restore the
# original protected state of the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the "as of" date in the confidentiality / disclosure footer text.
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-22 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) figures for each holding.
$ws.Range("D2").Value = 0.4897182876556266
$ws.Range("E2").Value = 0.001980982567353573

$ws.Range("D3").Value = 0.3329463477088175
$ws.Range("E3").Value = -0.005980861244019087

$ws.Range("D4").Value = 0.09287399692751711
$ws.Range("E4").Value = 0.005580142389840281

$ws.Range("D5").Value = 0.05492685461105393
$ws.Range("E5").Value = 0.003602138043225622

$ws.Range("D6").Value = 0.0295345130969849
$ws.Range("E6").Value = -0.004213938411669504

$ws.Range("E7").Value = -0.0004295348960403622

# Restore worksheet protection.
$ws.Protect()
